$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared text blocks reused from existing rows
$scaledSpeed = "scaled speed`nweekday o.h.`ndaypart o.h."
$lstmLayers  = "lstm(50)+do(.3)`nlstm/50)+do(.3)`nlstm/33)"
$oneWeek     = "1 week"
$firstWeekJune = "First 7 days of June"
$oneHour     = "1h back`n1h forward"

# New comment / timestamp text introduced for models 3 and 4
$febMarAprMay       = "Feb March April May"
$febMarAprMayOctNov = "Feb March April May Oct Nov"
$comment471         = "Increasing the training data hasn't done much on improving the accuracy. The increase is so slight that it can be ignored"
$comment1745a       = "Same data set gave a greater loss than 471 for 1745. "
$comment1745b       = "Increasing training data for 1745 has improved its performance some"

# Row 4 - Model 3, Speed ID 471
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 471
$ws.Range("C4").Value = $scaledSpeed
$ws.Range("D4").Value = $lstmLayers
$ws.Range("E4").Value = 50
$ws.Range("F4").Value = $oneHour
$ws.Range("G4").Value = $oneWeek
$ws.Range("H4").Value = $febMarAprMay
$ws.Range("I4").Value = $firstWeekJune
$ws.Range("J4").Value = 12.95
$ws.Range("K4").Value = 16.012
$ws.Range("L4").Value = 24
$ws.Range("M4").Value = $comment471

# Row 5 - Model 3, Speed ID 1745
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1745
$ws.Range("C5").Value = $scaledSpeed
$ws.Range("D5").Value = $lstmLayers
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = $oneHour
$ws.Range("G5").Value = $oneWeek
$ws.Range("H5").Value = $febMarAprMay
$ws.Range("I5").Value = $firstWeekJune
$ws.Range("J5").Value = 22.98
$ws.Range("K5").Value = 27.8
$ws.Range("L5").Value = 22.4
$ws.Range("M5").Value = $comment1745a

# Row 6 - Model 4, Speed ID 1745
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1745
$ws.Range("C6").Value = $scaledSpeed
$ws.Range("D6").Value = $lstmLayers
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = $oneHour
$ws.Range("G6").Value = $oneWeek
$ws.Range("H6").Value = $febMarAprMayOctNov
$ws.Range("I6").Value = $firstWeekJune
$ws.Range("J6").Value = 23.99
$ws.Range("K6").Value = 23.67
$ws.Range("L6").Value = 27.6
$ws.Range("M6").Value = $comment1745b

# Keep selection consistent with the author's saved view
$ws.Range("M6").Select()
